# Auto-generated edit script for unitTest_ws.xlsx (#system sheet reference tables)
# Reflects: removal of duplicate "clear(variables)" from base; new commands
# terminate/assertPath/assertAttributeContain/saveSelectedText/saveSelectedValue
# inserted (alphabetically) into external/io/web lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- Update "#system" sheet reference-table cell values ---
$ws.Cells.Item(18,6).Value = "clear(vars)"  # F18
$ws.Cells.Item(19,6).Value = "failImmediate(text)"  # F19
$ws.Cells.Item(20,6).Value = "incrementChar(var,amount,config)"  # F20
$ws.Cells.Item(21,6).Value = "macro(file,sheet,name)"  # F21
$ws.Cells.Item(22,6).Value = "outputToCloud(resource)"  # F22
$ws.Cells.Item(23,6).Value = "prependText(var,prependWith)"  # F23
$ws.Cells.Item(24,6).Value = "repeatUntil(steps,maxWaitMs)"  # F24
$ws.Cells.Item(25,6).Value = "save(var,value)"  # F25
$ws.Cells.Item(26,6).Value = "saveCount(text,regex,saveVar)"  # F26
$ws.Cells.Item(27,6).Value = "saveMatches(text,regex,saveVar)"  # F27
$ws.Cells.Item(28,6).Value = "saveReplace(text,regex,replace,saveVar)"  # F28
$ws.Cells.Item(29,6).Value = "saveVariablesByPrefix(var,prefix)"  # F29
$ws.Cells.Item(30,6).Value = "saveVariablesByRegex(var,regex)"  # F30
$ws.Cells.Item(31,6).Value = "section(steps)"  # F31
$ws.Cells.Item(32,6).Value = "split(text,delim,saveVar)"  # F32
$ws.Cells.Item(33,6).Value = "startRecording()"  # F33
$ws.Cells.Item(34,6).Value = "stopRecording()"  # F34
$ws.Cells.Item(35,6).Value = "substringAfter(text,delim,saveVar)"  # F35
$ws.Cells.Item(36,6).Value = "substringBefore(text,delim,saveVar)"  # F36
$ws.Cells.Item(37,6).Value = "substringBetween(text,start,end,saveVar)"  # F37
$ws.Cells.Item(38,6).Value = "verbose(text)"  # F38
$ws.Cells.Item(39,6).Value = "waitFor(waitMs)"  # F39
$ws.Cells.Item(40,6).ClearContents()  # F40 (was "waitFor(waitMs)")
$ws.Cells.Item(6,10).Value = "terminate(programName)"  # J6
$ws.Cells.Item(4,12).Value = "assertPath(path)"  # L4
$ws.Cells.Item(5,12).Value = "assertReadableFile(file,minByte)"  # L5
$ws.Cells.Item(6,12).Value = "base64(var,file)"  # L6
$ws.Cells.Item(7,12).Value = "compare(expected,actual,failFast)"  # L7
$ws.Cells.Item(8,12).Value = "copyFiles(source,target)"  # L8
$ws.Cells.Item(9,12).Value = "copyFilesByRegex(sourceDir,regex,target)"  # L9
$ws.Cells.Item(10,12).Value = "count(var,path,pattern)"  # L10
$ws.Cells.Item(11,12).Value = "deleteFiles(location,recursive)"  # L11
$ws.Cells.Item(12,12).Value = "deleteFilesByRegex(sourceDir,regex)"  # L12
$ws.Cells.Item(13,12).Value = "filter(source,target,matchPattern)"  # L13
$ws.Cells.Item(14,12).Value = "makeDirectory(source)"  # L14
$ws.Cells.Item(15,12).Value = "moveFiles(source,target)"  # L15
$ws.Cells.Item(16,12).Value = "moveFilesByRegex(sourceDir,regex,target)"  # L16
$ws.Cells.Item(17,12).Value = "readFile(var,file)"  # L17
$ws.Cells.Item(18,12).Value = "readProperty(var,file,property)"  # L18
$ws.Cells.Item(19,12).Value = "rename(target,newName)"  # L19
$ws.Cells.Item(20,12).Value = "saveDiff(var,expected,actual)"  # L20
$ws.Cells.Item(21,12).Value = "saveFileMeta(var,file)"  # L21
$ws.Cells.Item(22,12).Value = "saveMatches(var,path,fileFilter,textFilter)"  # L22
$ws.Cells.Item(23,12).Value = "searchAndReplace(file,config,saveAs)"  # L23
$ws.Cells.Item(24,12).Value = "unzip(zipFile,target)"  # L24
$ws.Cells.Item(25,12).Value = "validate(var,profile,inputFile)"  # L25
$ws.Cells.Item(26,12).Value = "writeBase64decode(encodedSource,decodedTarget,append)"  # L26
$ws.Cells.Item(27,12).Value = "writeFile(file,content,append)"  # L27
$ws.Cells.Item(28,12).Value = "writeFileAsIs(file,content,append)"  # L28
$ws.Cells.Item(29,12).Value = "writeProperty(file,property,value)"  # L29
$ws.Cells.Item(30,12).Value = "zip(filePattern,zipFile)"  # L30
$ws.Cells.Item(99,26).Value = "saveSelectedText(var,locator)"  # Z99
$ws.Cells.Item(100,26).Value = "saveSelectedValue(var,locator)"  # Z100
$ws.Cells.Item(101,26).Value = "saveTableAsCsv(locator,nextPageLocator,file)"  # Z101
$ws.Cells.Item(102,26).Value = "saveText(var,locator)"  # Z102
$ws.Cells.Item(103,26).Value = "saveTextArray(var,locator)"  # Z103
$ws.Cells.Item(104,26).Value = "saveTextSubstringAfter(var,locator,delim)"  # Z104
$ws.Cells.Item(105,26).Value = "saveTextSubstringBefore(var,locator,delim)"  # Z105
$ws.Cells.Item(106,26).Value = "saveTextSubstringBetween(var,locator,start,end)"  # Z106
$ws.Cells.Item(107,26).Value = "saveValue(var,locator)"  # Z107
$ws.Cells.Item(108,26).Value = "saveValues(var,locator)"  # Z108
$ws.Cells.Item(109,26).Value = "screenshot(file,locator)"  # Z109
$ws.Cells.Item(110,26).Value = "scrollElement(locator,xOffset,yOffset)"  # Z110
$ws.Cells.Item(111,26).Value = "scrollLeft(locator,pixel)"  # Z111
$ws.Cells.Item(112,26).Value = "scrollPage(xOffset,yOffset)"  # Z112
$ws.Cells.Item(113,26).Value = "scrollRight(locator,pixel)"  # Z113
$ws.Cells.Item(114,26).Value = "scrollTo(locator)"  # Z114
$ws.Cells.Item(115,26).Value = "select(locator,text)"  # Z115
$ws.Cells.Item(116,26).Value = "selectFrame(locator)"  # Z116
$ws.Cells.Item(117,26).Value = "selectMulti(locator,array)"  # Z117
$ws.Cells.Item(118,26).Value = "selectMultiOptions(locator)"  # Z118
$ws.Cells.Item(119,26).Value = "selectText(locator)"  # Z119
$ws.Cells.Item(120,26).Value = "selectWindow(winId)"  # Z120
$ws.Cells.Item(121,26).Value = "selectWindowAndWait(winId,waitMs)"  # Z121
$ws.Cells.Item(122,26).Value = "selectWindowByIndex(index)"  # Z122
$ws.Cells.Item(123,26).Value = "selectWindowByIndexAndWait(index,waitMs)"  # Z123
$ws.Cells.Item(124,26).Value = "toggleSelections(locator)"  # Z124
$ws.Cells.Item(125,26).Value = "type(locator,value)"  # Z125
$ws.Cells.Item(126,26).Value = "typeKeys(locator,value)"  # Z126
$ws.Cells.Item(127,26).Value = "uncheckAll(locator)"  # Z127
$ws.Cells.Item(128,26).Value = "unselectAllText()"  # Z128
$ws.Cells.Item(129,26).Value = "updateAttribute(locator,attrName,value)"  # Z129
$ws.Cells.Item(130,26).Value = "upload(fieldLocator,file)"  # Z130
$ws.Cells.Item(131,26).Value = "verifyContainText(locator,text)"  # Z131
$ws.Cells.Item(132,26).Value = "verifyText(locator,text)"  # Z132
$ws.Cells.Item(133,26).Value = "wait(waitMs)"  # Z133
$ws.Cells.Item(134,26).Value = "waitForElementPresent(locator)"  # Z134
$ws.Cells.Item(135,26).Value = "waitForPopUp(winId,waitMs)"  # Z135
$ws.Cells.Item(136,26).Value = "waitForTextPresent(text)"  # Z136
$ws.Cells.Item(137,26).Value = "waitForTitle(text)"  # Z137

# --- Update defined names to reflect new range extents ---
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$39"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$6"
$wb.Names.Item("io").RefersTo = "='#system'!`$L`$2:`$L`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$137"
